# Updated cryptos list on Wed Nov  6 10:58:46 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "73.840.56"
$ws.Range("E2").Value = "  +7.44%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.625.63"
$ws.Range("E3").Value = "  +7.81%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "185.22"
$ws.Range("E5").Value = "  +14.50%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "582.76"
$ws.Range("E6").Value = "  +4.32%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("E8").Value = "  +4.33%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.201"
$ws.Range("E9").Value = "  +19.43%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.623.80"
$ws.Range("E10").Value = "  +7.78%  "

$ws.Range("E11").Value = "  +0.38%  "

$ws.Range("E12").Value = "  +8.21%  "

$ws.Range("E13").Value = "  +3.86%  "

$ws.Range("E14").Value = "  +6.93%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.107.46"
$ws.Range("E15").Value = "  +7.73%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "73.703.73"
$ws.Range("E16").Value = "  +7.39%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.19"
$ws.Range("E17").Value = "  +12.85%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.626.16"
$ws.Range("E18").Value = "  +7.76%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.05"
$ws.Range("E19").Value = "  +30.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.84"
$ws.Range("E20").Value = "  +12.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "371.47"
$ws.Range("E21").Value = "  +9.53%  "

$ws.Range("E22").Value = "  +19.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.09"
$ws.Range("E23").Value = "  +6.83%  "

$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.84"
$ws.Range("E25").Value = "  +4.45%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.14"
$ws.Range("E26").Value = "  +11.85%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.38"
$ws.Range("E27").Value = "  +14.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.762.07"
$ws.Range("E28").Value = "  +7.69%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.01"
$ws.Range("E29").Value = "  -0.03%  "

$ws.Range("E30").Value = "  +15.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "522.72"
$ws.Range("E31").Value = "  +22.20%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.39"
$ws.Range("E32").Value = "  +19.83%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.65"
$ws.Range("E33").Value = "  +7.42%  "

$ws.Range("E34").Value = "  +9.17%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.119"
$ws.Range("E36").Value = "  +13.31%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.09"
$ws.Range("E37").Value = "  +1.39%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.16"
$ws.Range("E38").Value = "  +6.61%  "

$ws.Range("E39").Value = "  +1.47%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.91"
$ws.Range("E41").Value = "  +13.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.330"
$ws.Range("E42").Value = "  +10.20%  "

$ws.Range("E43").Value = "  +10.84%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "161.33"
$ws.Range("E44").Value = "  +23.19%  "

$ws.Range("E45").Value = "  +11.24%  "

$ws.Range("E46").Value = "  +15.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.96"
$ws.Range("E47").Value = "  +4.03%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0852"
$ws.Range("E48").Value = "  +18.38%  "

$ws.Range("E50").Value = "  +10.33%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.71"
$ws.Range("E51").Value = "  +22.85%  "
